$d = $word.ActiveDocument

function Get-ParaByText($marker) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$marker*") {
            return $p
        }
    }
    return $null
}

# 1) Remove stray "Pasted image 20241225095402.png" paragraph (style BodyText).
#    The following paragraph is also BodyText, so a plain merge-delete keeps the
#    correct style on the surviving paragraph.
$p = Get-ParaByText "Pasted image 20241225095402.png"
$p.Range.Delete()

# 2) Remove stray "Pasted image 20241225095443.png" paragraph (style FirstParagraph).
#    The following paragraph ("Các cấp độ WBS ...") is BodyText, but after the
#    merge the diff keeps it styled as FirstParagraph (the deleted paragraph's
#    own style), so we restyle the surviving merged paragraph explicitly.
$p = Get-ParaByText "Pasted image 20241225095443.png"
$p.Range.Delete()
$merged = Get-ParaByText "Các cấp độ WBS là yếu tố"
$merged.Style = "FirstParagraph"

# 3) Remove stray "Pasted image 20241225094201.png" paragraph (style BodyText).
#    The following paragraph is also BodyText.
$p = Get-ParaByText "Pasted image 20241225094201.png"
$p.Range.Delete()

# 4) Remove the paragraph containing both
#    "Pasted image 20241225094520.png" and "Pasted image 20241225094506.png"
#    (separated by a line break), which is immediately followed by bookmark
#    markers rather than another text paragraph.
$p = Get-ParaByText "Pasted image 20241225094520.png"
$p.Range.Delete()

# 5) Remove stray "Pasted image 20241225101324.png" paragraph (style
#    FirstParagraph), immediately followed by a bookmark start rather than
#    another text paragraph.
$p = Get-ParaByText "Pasted image 20241225101324.png"
$p.Range.Delete()

$d.Save()
